$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Highlight the blocked/404 issue row (row 7) in bold red, matching the
#    "pending attention" callout added by the author.
$rng = $ws.Range("B7:D7")
$rng.Font.Bold = $true
$rng.Font.Color = 255

# 2. Fill in the "plan type" (column D) for rows that previously had no
#    classification.
$ws.Range("D12").Value = "待续"
$ws.Range("D13").Value = "待续"
$ws.Range("D14").Value = "待续"
$ws.Range("D15").Value = "待续"
$ws.Range("D16").Value = "进行中"
$ws.Range("D17").Value = "待续"
$ws.Range("D18").Value = "暂停"
$ws.Range("D19").Value = "暂停"
$ws.Range("D20").Value = "暂停"

# 3. Add three new feature rows (sso module lookup/register + international),
#    duplicating the formatting of the last existing data row.
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(21).Insert(-4121)
$ws.Rows.Item(21).RowHeight = 24.75

$ws.Rows.Item(20).Copy()
$ws.Rows.Item(22).Insert(-4121)
$ws.Rows.Item(22).RowHeight = 24.75

$ws.Rows.Item(20).Copy()
$ws.Rows.Item(23).Insert(-4121)
$ws.Rows.Item(23).RowHeight = 24.75

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "icustom.sso"
$ws.Range("D21").Value = "待续"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "补齐数据模型和设计"
$ws.Range("D22").Value = "加急优先"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "补齐api文档"
$ws.Range("D23").Value = "加急优先"

# 4. Grow the table (and its autofilter) to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I23"))

# 5. Extend the "completion status" dropdown validation down to row 23
#    (re-created in the original E-then-D order).
$ws.Cells.Validation.Delete()
$ws.Range("E2:E23").Validation.Add(3, 1, 1, "=`$N`$3:`$N`$5")
$ws.Range("D1:D1048576").Validation.Add(3, 1, 1, "=`$L`$3:`$L`$11")

# 6. Move the active selection, matching the author's final cursor position.
$ws.Range("B6").Select()
